$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Helper: apply one of the "new" cell styles (matching the target
# OOXML cellXfs entries 19-33) to a given cell address.
# ---------------------------------------------------------------
function Set-NewStyle {
    param($addr, $styleId)

    $c = $ws.Range($addr)
    $c.Font.Name = "Arial"

    switch ($styleId) {
        19 { # fontId5 bold sz10 FF002060, halign left, valign center
            $c.Font.Bold = $true
            $c.Font.Size = 10
            $c.Font.Color = 6299648
            $c.HorizontalAlignment = -4131
            $c.VerticalAlignment = -4108
        }
        20 { # fontId4 bold sz9 theme4 tint-0.5, halign right, valign center
            $c.Font.Bold = $true
            $c.Font.Size = 9
            $c.Font.Color = 6567712
            $c.HorizontalAlignment = -4152
            $c.VerticalAlignment = -4108
        }
        21 { # fontId6 bold sz9 black, halign right, valign center
            $c.Font.Bold = $true
            $c.Font.Size = 9
            $c.Font.Color = 0
            $c.HorizontalAlignment = -4152
            $c.VerticalAlignment = -4108
        }
        22 { # fontId6 bold sz9 black, halign center, valign center
            $c.Font.Bold = $true
            $c.Font.Size = 9
            $c.Font.Color = 0
            $c.HorizontalAlignment = -4108
            $c.VerticalAlignment = -4108
        }
        23 { # fontId1 regular sz9 black, halign right, valign center
            $c.Font.Bold = $false
            $c.Font.Size = 9
            $c.Font.Color = 0
            $c.HorizontalAlignment = -4152
            $c.VerticalAlignment = -4108
        }
        24 { # fontId1 regular sz9 black, numFmt 0.00, halign right, valign center
            $c.Font.Bold = $false
            $c.Font.Size = 9
            $c.Font.Color = 0
            $c.HorizontalAlignment = -4152
            $c.VerticalAlignment = -4108
            $c.NumberFormat = "0.00"
        }
        25 { # fontId1 regular sz9 black, thin bottom border FF002060, halign right, valign center
            $c.Font.Bold = $false
            $c.Font.Size = 9
            $c.Font.Color = 0
            $c.HorizontalAlignment = -4152
            $c.VerticalAlignment = -4108
            $c.Borders.Item(9).LineStyle = 1
            $c.Borders.Item(9).Weight = 2
            $c.Borders.Item(9).Color = 6299648
        }
        26 { # fontId4 bold sz9 theme4 tint-0.5, numFmt 43 (Millares), halign right, valign center
            $c.Font.Bold = $true
            $c.Font.Size = 9
            $c.Font.Color = 6567712
            $c.HorizontalAlignment = -4152
            $c.VerticalAlignment = -4108
            $c.NumberFormat = '_-* #,##0.00_-;\-* #,##0.00_-;_-* "-"??_-;_-@_-'
        }
        27 { # fontId6 bold sz9 black, numFmt 43, halign right, valign center
            $c.Font.Bold = $true
            $c.Font.Size = 9
            $c.Font.Color = 0
            $c.HorizontalAlignment = -4152
            $c.VerticalAlignment = -4108
            $c.NumberFormat = '_-* #,##0.00_-;\-* #,##0.00_-;_-* "-"??_-;_-@_-'
        }
        28 { # fontId1 regular sz9 black, numFmt 43, halign center, valign center
            $c.Font.Bold = $false
            $c.Font.Size = 9
            $c.Font.Color = 0
            $c.HorizontalAlignment = -4108
            $c.VerticalAlignment = -4108
            $c.NumberFormat = '_-* #,##0.00_-;\-* #,##0.00_-;_-* "-"??_-;_-@_-'
        }
        29 { # fontId4 bold sz9 theme4 tint-0.5, numFmt 43, medium bottom border FF002060, halign right, valign center
            $c.Font.Bold = $true
            $c.Font.Size = 9
            $c.Font.Color = 6567712
            $c.HorizontalAlignment = -4152
            $c.VerticalAlignment = -4108
            $c.NumberFormat = '_-* #,##0.00_-;\-* #,##0.00_-;_-* "-"??_-;_-@_-'
            $c.Borders.Item(9).LineStyle = 1
            $c.Borders.Item(9).Weight = -4138
            $c.Borders.Item(9).Color = 6299648
        }
        30 { # fontId1 regular sz9 black, numFmt 43, halign right, valign center
            $c.Font.Bold = $false
            $c.Font.Size = 9
            $c.Font.Color = 0
            $c.HorizontalAlignment = -4152
            $c.VerticalAlignment = -4108
            $c.NumberFormat = '_-* #,##0.00_-;\-* #,##0.00_-;_-* "-"??_-;_-@_-'
        }
        31 { # fontId1 regular sz9 black, numFmt 43, thin bottom border FF002060, halign right, valign center
            $c.Font.Bold = $false
            $c.Font.Size = 9
            $c.Font.Color = 0
            $c.HorizontalAlignment = -4152
            $c.VerticalAlignment = -4108
            $c.NumberFormat = '_-* #,##0.00_-;\-* #,##0.00_-;_-* "-"??_-;_-@_-'
            $c.Borders.Item(9).LineStyle = 1
            $c.Borders.Item(9).Weight = 2
            $c.Borders.Item(9).Color = 6299648
        }
        32 { # fontId1 regular sz9 black, thin bottom border default/automatic, halign right, valign center
            $c.Font.Bold = $false
            $c.Font.Size = 9
            $c.Font.Color = 0
            $c.HorizontalAlignment = -4152
            $c.VerticalAlignment = -4108
            $c.Borders.Item(9).LineStyle = 1
            $c.Borders.Item(9).Weight = 2
            $c.Borders.Item(9).ColorIndex = -4105
        }
        33 { # fontId1 regular sz9 black, numFmt 0.00, thin bottom border default/automatic, halign right, valign center
            $c.Font.Bold = $false
            $c.Font.Size = 9
            $c.Font.Color = 0
            $c.HorizontalAlignment = -4152
            $c.VerticalAlignment = -4108
            $c.NumberFormat = "0.00"
            $c.Borders.Item(9).LineStyle = 1
            $c.Borders.Item(9).Weight = 2
            $c.Borders.Item(9).ColorIndex = -4105
        }
        5 { # fontId1 regular sz9 black, halign center, valign center (plain, blank)
            $c.Font.Bold = $false
            $c.Font.Size = 9
            $c.Font.Color = 0
            $c.HorizontalAlignment = -4108
            $c.VerticalAlignment = -4108
        }
        6 { # fontId4 bold sz9 theme4 tint-0.5, medium bottom border FF002060, halign right, valign center (existing style)
            $c.Font.Bold = $true
            $c.Font.Size = 9
            $c.Font.Color = 6567712
            $c.HorizontalAlignment = -4152
            $c.VerticalAlignment = -4108
            $c.Borders.Item(9).LineStyle = 1
            $c.Borders.Item(9).Weight = -4138
            $c.Borders.Item(9).Color = 6299648
        }
    }
}

# ---------------------------------------------------------------
# Table "2 - Frequency distribution table" (rows 12-24)
# ---------------------------------------------------------------
$ws.Range("D12").Value = "2 - Frequency distribution table"
Set-NewStyle "D12" 19

Set-NewStyle "D13" 5

$ws.Range("D14").Value = "Desired intervals"
Set-NewStyle "D14" 20
$ws.Range("E14").Value = 6
Set-NewStyle "E14" 21
Set-NewStyle "F14" 22

Set-NewStyle "D16" 5

$ws.Range("E15").Formula = '=ROUNDUP((B32-B13)/E14,0)'
Set-NewStyle "E15" 21

$ws.Range("D17").Value = "Interval start"
$ws.Range("E17").Value = "Interval end"
$ws.Range("F17").Value = "Frequency"
$ws.Range("G17").Value = "Relative frequency"
Set-NewStyle "D17" 6
Set-NewStyle "E17" 6
Set-NewStyle "F17" 6
Set-NewStyle "G17" 6

$ws.Range("D18").Value = 8
$ws.Range("E18").Formula = '=D18+E15'
$ws.Range("F18").Formula = '=COUNTIFS($B$13:$B$32,">="&D18,$B$13:$B$32,"<="&E18)'
$ws.Range("G18").Formula = '=F18/COUNT($B$13:$B$32)'
Set-NewStyle "D18" 23
Set-NewStyle "E18" 23
Set-NewStyle "F18" 23
Set-NewStyle "G18" 24

$ws.Range("D19").Formula = '=E18'
$ws.Range("E19").Formula = '=D19+E15'
Set-NewStyle "D19" 23
Set-NewStyle "E19" 23

$ws.Range("D20").Formula = '=E19'
$ws.Range("E20").Formula = '=D20+E15'
Set-NewStyle "D20" 23
Set-NewStyle "E20" 23

$ws.Range("D21").Formula = '=E20'
$ws.Range("E21").Formula = '=D21+E15'
Set-NewStyle "D21" 23
Set-NewStyle "E21" 23

$ws.Range("D22").Formula = '=E21'
$ws.Range("E22").Formula = '=D22+E15'
Set-NewStyle "D22" 23
Set-NewStyle "E22" 23

$ws.Range("D23").Formula = '=E22'
$ws.Range("E23").Formula = '=D23+E15'
Set-NewStyle "D23" 25
Set-NewStyle "E23" 25

# Shared formulas F19:F23 and G19:G24
$ws.Range("F19:F23").Formula = '=COUNTIFS($B$13:$B$32,">="&D19,$B$13:$B$32,"<="&E19)'
$ws.Range("G19:G24").Formula = '=F19/COUNT($B$13:$B$32)'

Set-NewStyle "F19" 23
Set-NewStyle "F20" 23
Set-NewStyle "F21" 23
Set-NewStyle "F22" 23
Set-NewStyle "F23" 32
Set-NewStyle "G19" 24
Set-NewStyle "G20" 24
Set-NewStyle "G21" 24
Set-NewStyle "G22" 24
Set-NewStyle "G23" 33
Set-NewStyle "G24" 30

$ws.Range("F24").Formula = '=SUM(F18:F23)'
Set-NewStyle "F24" 23

# ---------------------------------------------------------------
# Table "3 - Frequency distribution table" (rows 26-38)
# ---------------------------------------------------------------
$ws.Range("D26").Value = "3 - Frequency distribution table"
Set-NewStyle "D26" 19

Set-NewStyle "D27" 5

$ws.Range("D28").Value = "Desired intervals"
Set-NewStyle "D28" 20
$ws.Range("E28").Value = 6
Set-NewStyle "E28" 21
Set-NewStyle "F28" 22

$ws.Range("D29").Value = "Interval width"
Set-NewStyle "D29" 26
$ws.Range("E29").Formula = '=(B32-B13)/E28'
Set-NewStyle "E29" 27

Set-NewStyle "D30" 28
Set-NewStyle "E30" 28

$ws.Range("D31").Value = "Interval start"
$ws.Range("E31").Value = "Interval end"
$ws.Range("F31").Value = "Frequency"
$ws.Range("G31").Value = "Relative frequency"
Set-NewStyle "D31" 29
Set-NewStyle "E31" 29
Set-NewStyle "F31" 6
Set-NewStyle "G31" 6

$ws.Range("D32").Value = 8
$ws.Range("E32").Formula = '=D32+E29'
$ws.Range("F32").Formula = '=COUNTIFS($B$13:$B$32,">="&D32,$B$13:$B$32,"<="&E32)'
$ws.Range("G32").Formula = '=F32/COUNT($B$13:$B$32)'
Set-NewStyle "D32" 30
Set-NewStyle "E32" 30
Set-NewStyle "F32" 23
Set-NewStyle "G32" 24

$ws.Range("D33").Formula = '=E32'
$ws.Range("E33").Formula = '=D33+E29'
Set-NewStyle "D33" 30
Set-NewStyle "E33" 30

$ws.Range("D34").Formula = '=E33'
$ws.Range("E34").Formula = '=D34+E29'
Set-NewStyle "D34" 30
Set-NewStyle "E34" 30

$ws.Range("D35").Formula = '=E34'
$ws.Range("E35").Formula = '=D35+E29'
Set-NewStyle "D35" 30
Set-NewStyle "E35" 30

$ws.Range("D36").Formula = '=E35'
$ws.Range("E36").Formula = '=D36+E29'
Set-NewStyle "D36" 30
Set-NewStyle "E36" 30

$ws.Range("D37").Formula = '=E36'
$ws.Range("E37").Formula = '=D37+E29'
Set-NewStyle "D37" 31
Set-NewStyle "E37" 31

# Shared formulas F33:F37 and G33:G38
$ws.Range("F33:F37").Formula = '=COUNTIFS($B$13:$B$32,">="&D33,$B$13:$B$32,"<="&E33)'
$ws.Range("G33:G38").Formula = '=F33/COUNT($B$13:$B$32)'

Set-NewStyle "F33" 23
Set-NewStyle "F34" 23
Set-NewStyle "F35" 23
Set-NewStyle "F36" 23
Set-NewStyle "F37" 32
Set-NewStyle "G33" 24
Set-NewStyle "G34" 24
Set-NewStyle "G35" 24
Set-NewStyle "G36" 24
Set-NewStyle "G37" 33
Set-NewStyle "G38" 30

$ws.Range("F38").Formula = '=SUM(F32:F37)'
Set-NewStyle "F38" 23

Write-Host "done"
